$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the name column (A) for data rows 2-5: replace old name with new name
$ws.Range("A2").Value = "Duvan Camilo"
$ws.Range("A3").Value = "Duvan Camilo"
$ws.Range("A4").Value = "Duvan Camilo"
$ws.Range("A5").Value = "Duvan Camilo"

# Update the country column (D) for data rows 2-5: all become Colombia
$ws.Range("D2").Value = "Colombia"
$ws.Range("D3").Value = "Colombia"
$ws.Range("D4").Value = "Colombia"
$ws.Range("D5").Value = "Colombia"

$ws.Range("C11").Select()
